# Add a new "october-2025" worksheet at the end of the workbook, based on the
# existing "september-2025" sheet (same layout/formatting), then update its
# single cell with the new Corporation Tax summary line.

$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore focus at the end.
$originalActiveSheetName = $wb.ActiveSheet.Name

# Copy the last monthly sheet ("september-2025") and place the copy right
# after it, i.e. at the very end of the workbook.
$src = $wb.Worksheets.Item("september-2025")
$src.Copy($null, $src)

# The newly created sheet is now the last worksheet in the workbook.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "october-2025"

# Replace the copied text with the October 2025 Corporation Tax figures.
$newSheet.Range("A1").Value = "Corporation Tax                                19,378           18,230                    1,148               6.3%"

# Restore the originally active sheet so we don't leave the new sheet selected.
$wb.Worksheets.Item($originalActiveSheetName).Activate()
